$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Presupuesto")

# Insert a new row before row 6, shifting existing rows 6-8 down to 7-9
$ws.Rows.Item(6).Insert()

# Populate the new row 6 with the "Interés por caución" entry
$ws.Range("A6").Value = "Interés por caución"
$ws.Range("B6").Value = "Columna Int. por caución en flujo por mes: cálculo de interés mensual por reinversión del sobrante a un día con tasa de serie de cauciones. Carga de Excel al refrescar, modal de marcha de cálculo (G/P acum, Base, Tasa, Int T). Incluye soporte para múltiples formatos de fecha y columna tasa_diaria."
$ws.Range("C6").Value = 50000
